$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24

# Plain text values - no auto-detection issues for these.
$ws.Cells.Item($row, 1).Value = "M1"
$ws.Cells.Item($row, 2).Value = "Marcas"
$ws.Cells.Item($row, 3).Value = "Bool"
$ws.Cells.Item($row, 4).Value = "%M2000.1"

# Empty-string text cells (Comment / Typeobject ID / Version ID columns).
# A bare "" assignment leaves the cell completely blank (no stored value),
# so force literal text entry via a leading quote-prefix, then strip the
# resulting "stored as text" cell formatting back off.
foreach ($col in 5, 9, 10) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'"
    $cell.ClearFormats()
}

# "True" text cells (Hmi Visible / Accessible / Writeable columns).
# Assigning the literal word True/TRUE gets auto-coerced to a Boolean, so
# instead compute it as a text formula in a scratch cell and paste back
# only the resulting value (which keeps it a genuine text cell).
$helper = $ws.Cells.Item(1, 300)
$helper.Formula = '="True"'
$helper.Copy()
foreach ($col in 6, 7, 8) {
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)
}
$helper.Clear()
$excel.CutCopyMode = $false
